$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.075.05"
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.892.40"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "306.51"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +0.09%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5221"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.44%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3757"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.67%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07263"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.02%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.09"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  +0.31%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08190"
$c.ClearFormats()
$ws.Range("E12").Value = "  +6.58%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.933.08"
$c.ClearFormats()
$ws.Range("E13").Value = "  +1.41%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "96.29"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.288"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  +0.12%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008580"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.26%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.58"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.90%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.101.90"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.58%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.082"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +0.56%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.403"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.45%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "148.41"
$c.ClearFormats()
$ws.Range("E24").Value = "  +1.82%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.285"
$c.ClearFormats()
$ws.Range("E25").Value = "  -1.09%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.16"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.08%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.729"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "114.90"
$c.ClearFormats()
$ws.Range("E28").Value = "  +0.12%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.777"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.834"
$c.ClearFormats()
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  +0.03%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.05033"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.29%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7874"
$c.ClearFormats()
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("E34").Value = "  -2.36%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.979"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.68%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.422"
$c.ClearFormats()
$ws.Range("E36").Value = "  +3.30%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.599"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.5724"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.75%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01980"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.36%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.073"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.18%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "9.044"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.66%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.547"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.11%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "116.25"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.56%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1514"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.01%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4854"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.47%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.12%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.09"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.81%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.623"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.13%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "38.09"
$c.ClearFormats()
$ws.Range("E49").Value = "  +1.42%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "63.51"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.33%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05928"
$c.ClearFormats()
$ws.Range("E51").Value = "  -0.03%  "
